$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the content of rows 42 and 43 (columns A-D)
$ws.Range("A42").Value = "c0a3f3ed23f04247d92740a9502f8b57"
$ws.Range("B42").Value = "unassigned"
$ws.Range("C42").Value = "unassigned"
$ws.Range("D42").Value = "unassigned"

$ws.Range("A43").Value = "307c55294ffe3b8aa46fce358d55590e"
$ws.Range("B43").Value = "Homo sapiens"
$ws.Range("C43").Value = "Human"
$ws.Range("D43").Value = "Human"

# Move the Station18 (J) value from row 43 to row 42
$ws.Range("J42").Value = 0
$ws.Range("J43").ClearContents()
